$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "330.57"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.25%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.70%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.702"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.38%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08396"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.10%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.815"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.93%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "2.001"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.19%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.477"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.29%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.34%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9228"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.07%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1282"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.73%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.13%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09528"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.01%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "3.44%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.83%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001302"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.51%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.18%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.76%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.69%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.919"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.09%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.95%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2510"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.50%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04414"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.47%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001274"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.93%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004356"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.54%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.02%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02840"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.23%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05517"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.52%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007952"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.33%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.86%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008979"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-9.73%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002061"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.00%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01168"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.95%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006939"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.44%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.14%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003463"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "14.79%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.21%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.14%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.14%"
